# Updated capital structure database
# Apply new values to rows 2 and 3 for the Oman Investments & Asset Management sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G" = 1.068493150684932
    "H" = 1.068493150684932
    "I" = 0.6061643835616438
    "J" = 0.6061643835616438
    "K" = -0.925
    "L" = -0.3167808219178083
    "M" = 3.12
    "N" = 0.09842271293375394
    "O" = -3.372972972972973
    "P" = 3.12
    "Q" = 0.09842271293375394
    "R" = -3.372972972972973
    "U" = 0.364
    "V" = 0.01148264984227129
    "W" = -0.01089517078916372
    "X" = 0.1122247688927248
    "Y" = -0.1231199396818885
    "Z" = 0.02254304022234231
    "AA" = 0.01366478807998147
    "AB" = 0.07081857879852553
    "AC" = -0.05715379071854406
    "AD" = 47.8
    "AF" = 47.8
    "AG" = 47.436
    "AH" = 0.6012578616352201
    "AI" = 0.3740219092331769
    "AJ" = 0.5994237767893248
    "AK" = 0.372233905646756
    "AL" = 2.54
    "AM" = 2.54
    "AN" = 26.85393258426966
    "AO" = 0.6968503937007874
    "AP" = 26.64943820224719
    "AQ" = 0.6968503937007874
}

foreach ($col in $updates.Keys) {
    $value = $updates[$col]
    $ws.Range($col + "2").Value = $value
    $ws.Range($col + "3").Value = $value
}
